$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "Populacao 150 / Geracao 500" and "Populacao 200 / Geracao 375" dataset
# rows (5 and 6) were stored in the wrong order - swap them so the dataset columns
# stay in the right sequence (commit: "os dataSets ficarem com as colunas seguidas").
# Columns C:G hold the shared configuration values and are identical on both rows,
# so only A, B, I and J actually need to move.

# --- read the current values first (Value2 -- Value's getter is unreliable here) ---
$a5 = $ws.Cells.Item(5, 1).Value2   # População ...
$b5 = $ws.Cells.Item(5, 2).Value2   # Geração ...
$i5 = $ws.Cells.Item(5, 9).Value2   # Average (number)
$j5 = $ws.Cells.Item(5, 10).Value2  # StdDev (numeric-looking text)

$a6 = $ws.Cells.Item(6, 1).Value2
$b6 = $ws.Cells.Item(6, 2).Value2
$i6 = $ws.Cells.Item(6, 9).Value2
$j6 = $ws.Cells.Item(6, 10).Value2

# --- write them back swapped ---
$ws.Cells.Item(5, 1).Value2 = $a6
$ws.Cells.Item(5, 2).Value2 = $b6
$ws.Cells.Item(5, 9).Value2 = $i6

$ws.Cells.Item(6, 1).Value2 = $a5
$ws.Cells.Item(6, 2).Value2 = $b5
$ws.Cells.Item(6, 9).Value2 = $i5

# J5/J6 ("5.18135117512797", "5.6718603649948935", ...) are stored as TEXT even
# though they look numeric, so force text formatting while writing them back or
# Excel would silently reinterpret them as numbers.
$ws.Cells.Item(5, 10).NumberFormat = "@"
$ws.Cells.Item(5, 10).Value2 = $j6
$ws.Cells.Item(5, 10).Style = "Normal"

$ws.Cells.Item(6, 10).NumberFormat = "@"
$ws.Cells.Item(6, 10).Value2 = $j5
$ws.Cells.Item(6, 10).Style = "Normal"

# Update the active cell selection left behind when the author last saved the file
$ws.Range("F15").Select()
